$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Periodo Mora" value from 2507 to 2508 for all affected rows.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
$ws.Range("E19").Value = "2508"

# Update the "Valor Mora" amount for ROQUE ALEXANDER PAJARO ACEVEDO (row 18).
$ws.Range("G18").Value = 1423500
